# Updates every cached "1/15/14" date-field display text to "4/17/16"
# across slides, slide layouts, the slide master, the notes master and the
# handout master, and rewrites the slide-1 subtitle from
# "Assignment Handin & Peer Grading" to "Peer Grading".

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "1/15/14") {
                $tr.Text = "4/17/16"
            }
        }
    }
}

# Slide master
Update-DateShapes $p.SlideMaster.Shapes

# Slide layouts
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master
Update-DateShapes $p.NotesMaster.Shapes

# Handout master
Update-DateShapes $p.HandoutMaster.Shapes

# Slides
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    Update-DateShapes $p.Slides.Item($si).Shapes
}

# Slide 1 title/subtitle text change: "Assignment Handin & Peer Grading"
# becomes "Peer Grading", split across two runs ("Peer " + "Grading") to
# mirror the authored edit.
$titleShape = $p.Slides.Item(1).Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange
$tr.Text = "Peer Grading"
$tr.Characters(1, 5).Text = "Peer "
$tr.Characters(6, 7).Text = "Grading"
